# Update odds data on the active worksheet to match the latest FlashScore
# refresh for the affected matches (rows 4 and 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (America De Cali - Millonarios)
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5

# Row 12 (Penarol - Maldonado)
$ws.Range("G12").Value = 1.38
$ws.Range("H12").Value = 4.2
$ws.Range("I12").Value = 9.5
$ws.Range("J12").Value = 1.95
$ws.Range("L12").Value = 8.5
$ws.Range("M12").Value = 1.08
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
$ws.Range("W12").Value = 5
$ws.Range("AG12").Value = 17
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 126
$ws.Range("AK12").Value = 81
$ws.Range("AZ12").Value = 251
$ws.Range("BA12").Value = 301
